# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# --- YDS: append newly logged per-play yardage to the four running lists ---
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 4 2 2 4 2 3 7 9 8 -1 4 11 5 4 9 3 5 5 6 0 1 0 5"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 2 4 1 7 2 -2 5 3 16 4 6 8 12 6 13 4 -1 6 10 7 3 15 2 8"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 17 6 11 31 7 22 10 3 13 3 24 6 1 16 9 7 6 9 25 14 11 23 -6"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 3 43 5 5 4 22 8 4 8 7 9 2 7 1 9 12 5 5 3 5 5 5 3 3 6 12 3"

# --- OFF: season totals after Week 15 + simulated Week 16 ---
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = 181
$offWs.Range("F2").Value = 49
$offWs.Range("G2").Value = 55
$offWs.Range("N2").Value = 14
$offWs.Range("B3").Value = 11
$offWs.Range("C3").Value = 143
$offWs.Range("D3").Value = 8
$offWs.Range("E3").Value = 20
$offWs.Range("F3").Value = 92
$offWs.Range("G3").Value = 40
$offWs.Range("H3").Value = 26
$offWs.Range("I3").Value = 49
$offWs.Range("J3").Value = 52
$offWs.Range("L3").Value = 255
$offWs.Range("M3").Value = 168
$offWs.Range("Q3").Value = 453

# --- DEF: season totals after Week 15 + simulated Week 16 ---
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 153
$defWs.Range("D2").Value = 10
$defWs.Range("E2").Value = 3
$defWs.Range("F2").Value = 45
$defWs.Range("G2").Value = 45
$defWs.Range("H2").Value = 4
$defWs.Range("J2").Value = 33
$defWs.Range("N2").Value = 17
$defWs.Range("O2").Value = 25
$defWs.Range("P2").Value = 11
$defWs.Range("C3").Value = 152
$defWs.Range("D3").Value = 6
$defWs.Range("E3").Value = 24
$defWs.Range("F3").Value = 96
$defWs.Range("G3").Value = 38
$defWs.Range("I3").Value = 45
$defWs.Range("J3").Value = 54
$defWs.Range("L3").Value = 275
$defWs.Range("M3").Value = 176
$defWs.Range("Q3").Value = 500

# --- ST: season totals + appended per-kick lists ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 73
$stWs.Range("D2").Value = 43
$stWs.Range("F2").Value = 708
$stWs.Range("G2").Value = 691
$stWs.Range("J2").Value = 254
$stWs.Range("K2").Value = 231
$stWs.Range("B3").Value = 32
$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 66 64 62 61"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 24 26 34 21"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 35 21 20 0 9"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 47 59 40"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 0 0 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0"

# --- TURNS: season totals ---
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("D3").Value = 6
$turnsWs.Range("E3").Value = 9

# --- PEN: season totals ---
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 9
$penWs.Range("D2").Value = 5
$penWs.Range("D4").Value = 10

